$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (Key/Value) - data shifts up by one row.
$ws.Rows.Item(1).Delete()

# Clear style on B1 (the old B2 had no explicit style / "s" attribute)
$ws.Range("B1").ClearFormats()
